$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.7169949999999999"
$ws.Range("H2").Value = [double]"2.150985"
$ws.Range("I2").Value = [double]"0.003651663653539308"
$ws.Range("J2").Value = [double]"0.003651663653539308"
$ws.Range("M2").Value = [double]"28.31444233333334"
$ws.Range("N2").Value = [double]"84.94332700000001"
$ws.Range("O2").Value = [double]"0.2747173016130739"
$ws.Range("P2").Value = [double]"0.2747173016130739"
$ws.Range("Q2").Value = [double]"20.30131358078833"
$ws.Range("R2").Value = [double]"182.711822227095"
$ws.Range("S2").Value = [double]"0.001003175185298858"
$ws.Range("T2").Value = [double]"0.001003175185298858"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.7169949999999999"
$ws.Range("H3").Value = [double]"2.150985"
$ws.Range("I3").Value = [double]"0.003651663653539308"
$ws.Range("J3").Value = [double]"0.003651663653539308"
$ws.Range("O3").Value = [double]"0.2090339131726295"
$ws.Range("P3").Value = [double]"0.2090339131726295"
$ws.Range("Q3").Value = [double]"15.44738170992167"
$ws.Range("R3").Value = [double]"139.026435389295"
$ws.Range("S3").Value = [double]"0.0007633215430895828"
$ws.Range("T3").Value = [double]"0.0007633215430895828"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.7169949999999999"
$ws.Range("H4").Value = [double]"2.150985"
$ws.Range("I4").Value = [double]"0.003651663653539308"
$ws.Range("J4").Value = [double]"0.003651663653539308"
$ws.Range("M4").Value = [double]"5.413469333333334"
$ws.Range("N4").Value = [double]"16.240408"
$ws.Range("O4").Value = [double]"0.0525235026743817"
$ws.Range("P4").Value = [double]"0.0525235026743817"
$ws.Range("Q4").Value = [double]"3.881430444653334"
$ws.Range("R4").Value = [double]"34.93287400188"
$ws.Range("S4").Value = [double]"0.0001917981656726143"
$ws.Range("T4").Value = [double]"0.0001917981656726143"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.7169949999999999"
$ws.Range("H5").Value = [double]"2.150985"
$ws.Range("I5").Value = [double]"0.003651663653539308"
$ws.Range("J5").Value = [double]"0.003651663653539308"
$ws.Range("M5").Value = [double]"47.79503400000001"
$ws.Range("N5").Value = [double]"143.385102"
$ws.Range("O5").Value = [double]"0.4637252825399149"
$ws.Range("P5").Value = [double]"0.4637252825399149"
$ws.Range("Q5").Value = [double]"34.26880040283"
$ws.Range("R5").Value = [double]"308.41920362547"
$ws.Range("S5").Value = [double]"0.001693368759478253"
$ws.Range("T5").Value = [double]"0.001693368759478253"
$ws.Range("I6").Value = [double]"0.9751961860217362"
$ws.Range("J6").Value = [double]"0.9751961860217361"
$ws.Range("M6").Value = [double]"28.31444233333334"
$ws.Range("N6").Value = [double]"84.94332700000001"
$ws.Range("O6").Value = [double]"0.2747173016130739"
$ws.Range("P6").Value = [double]"0.2747173016130739"
$ws.Range("Q6").Value = [double]"5421.573686291572"
$ws.Range("R6").Value = [double]"48794.16317662416"
$ws.Range("S6").Value = [double]"0.2679032647672527"
$ws.Range("T6").Value = [double]"0.2679032647672526"
$ws.Range("I7").Value = [double]"0.9751961860217362"
$ws.Range("J7").Value = [double]"0.9751961860217361"
$ws.Range("O7").Value = [double]"0.2090339131726295"
$ws.Range("P7").Value = [double]"0.2090339131726295"
$ws.Range("S7").Value = [double]"0.2038490748751471"
$ws.Range("T7").Value = [double]"0.203849074875147"
$ws.Range("I8").Value = [double]"0.9751961860217362"
$ws.Range("J8").Value = [double]"0.9751961860217361"
$ws.Range("M8").Value = [double]"5.413469333333334"
$ws.Range("N8").Value = [double]"16.240408"
$ws.Range("O8").Value = [double]"0.0525235026743817"
$ws.Range("P8").Value = [double]"0.0525235026743817"
$ws.Range("Q8").Value = [double]"1036.55662872069"
$ws.Range("R8").Value = [double]"9329.009658486209"
$ws.Range("S8").Value = [double]"0.05122071948455949"
$ws.Range("T8").Value = [double]"0.05122071948455949"
$ws.Range("I9").Value = [double]"0.9751961860217362"
$ws.Range("J9").Value = [double]"0.9751961860217361"
$ws.Range("M9").Value = [double]"47.79503400000001"
$ws.Range("N9").Value = [double]"143.385102"
$ws.Range("O9").Value = [double]"0.4637252825399149"
$ws.Range("P9").Value = [double]"0.4637252825399149"
$ws.Range("Q9").Value = [double]"9151.665274535728"
$ws.Range("R9").Value = [double]"82364.98747082155"
$ws.Range("S9").Value = [double]"0.452223126894777"
$ws.Range("T9").Value = [double]"0.452223126894777"
$ws.Range("G10").Value = [double]"4.138615666666666"
$ws.Range("H10").Value = [double]"12.415847"
$ws.Range("I10").Value = [double]"0.02107801645190694"
$ws.Range("J10").Value = [double]"0.02107801645190694"
$ws.Range("M10").Value = [double]"28.31444233333334"
$ws.Range("N10").Value = [double]"84.94332700000001"
$ws.Range("O10").Value = [double]"0.2747173016130739"
$ws.Range("P10").Value = [double]"0.2747173016130739"
$ws.Range("Q10").Value = [double]"117.1825946336632"
$ws.Range("R10").Value = [double]"1054.643351702969"
$ws.Range("S10").Value = [double]"0.005790495803023855"
$ws.Range("T10").Value = [double]"0.005790495803023854"
$ws.Range("G11").Value = [double]"4.138615666666666"
$ws.Range("H11").Value = [double]"12.415847"
$ws.Range("I11").Value = [double]"0.02107801645190694"
$ws.Range("J11").Value = [double]"0.02107801645190694"
$ws.Range("O11").Value = [double]"0.2090339131726295"
$ws.Range("P11").Value = [double]"0.2090339131726295"
$ws.Range("Q11").Value = [double]"89.16488393037878"
$ws.Range("R11").Value = [double]"802.483955373409"
$ws.Range("S11").Value = [double]"0.004406020260859172"
$ws.Range("T11").Value = [double]"0.004406020260859172"
$ws.Range("G12").Value = [double]"4.138615666666666"
$ws.Range("H12").Value = [double]"12.415847"
$ws.Range("I12").Value = [double]"0.02107801645190694"
$ws.Range("J12").Value = [double]"0.02107801645190694"
$ws.Range("M12").Value = [double]"5.413469333333334"
$ws.Range("N12").Value = [double]"16.240408"
$ws.Range("O12").Value = [double]"0.0525235026743817"
$ws.Range("P12").Value = [double]"0.0525235026743817"
$ws.Range("Q12").Value = [double]"22.40426899395289"
$ws.Range("R12").Value = [double]"201.638420945576"
$ws.Range("S12").Value = [double]"0.001107091253482396"
$ws.Range("T12").Value = [double]"0.001107091253482396"
$ws.Range("G13").Value = [double]"4.138615666666666"
$ws.Range("H13").Value = [double]"12.415847"
$ws.Range("I13").Value = [double]"0.02107801645190694"
$ws.Range("J13").Value = [double]"0.02107801645190694"
$ws.Range("M13").Value = [double]"47.79503400000001"
$ws.Range("N13").Value = [double]"143.385102"
$ws.Range("O13").Value = [double]"0.4637252825399149"
$ws.Range("P13").Value = [double]"0.4637252825399149"
$ws.Range("Q13").Value = [double]"197.805276501266"
$ws.Range("R13").Value = [double]"1780.247488511394"
$ws.Range("S13").Value = [double]"0.009774409134541522"
$ws.Range("T13").Value = [double]"0.00977440913454152"
$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.014556"
$ws.Range("H14").Value = [double]"0.043668"
$ws.Range("I14").Value = [double]"7.413387281768795E-05"
$ws.Range("J14").Value = [double]"7.413387281768795E-05"
$ws.Range("M14").Value = [double]"28.31444233333334"
$ws.Range("N14").Value = [double]"84.94332700000001"
$ws.Range("O14").Value = [double]"0.2747173016130739"
$ws.Range("P14").Value = [double]"0.2747173016130739"
$ws.Range("Q14").Value = [double]"0.412145022604"
$ws.Range("R14").Value = [double]"3.709305203436"
$ws.Range("S14").Value = [double]"2.036585749860204E-05"
$ws.Range("T14").Value = [double]"2.036585749860204E-05"
$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.014556"
$ws.Range("H15").Value = [double]"0.043668"
$ws.Range("I15").Value = [double]"7.413387281768795E-05"
$ws.Range("J15").Value = [double]"7.413387281768795E-05"
$ws.Range("O15").Value = [double]"0.2090339131726295"
$ws.Range("P15").Value = [double]"0.2090339131726295"
$ws.Range("Q15").Value = [double]"0.313603425644"
$ws.Range("R15").Value = [double]"2.822430830796"
$ws.Range("S15").Value = [double]"1.549649353372334E-05"
$ws.Range("T15").Value = [double]"1.549649353372334E-05"
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.014556"
$ws.Range("H16").Value = [double]"0.043668"
$ws.Range("I16").Value = [double]"7.413387281768795E-05"
$ws.Range("J16").Value = [double]"7.413387281768795E-05"
$ws.Range("M16").Value = [double]"5.413469333333334"
$ws.Range("N16").Value = [double]"16.240408"
$ws.Range("O16").Value = [double]"0.0525235026743817"
$ws.Range("P16").Value = [double]"0.0525235026743817"
$ws.Range("Q16").Value = [double]"0.07879845961600002"
$ws.Range("R16").Value = [double]"0.7091861365440001"
$ws.Range("S16").Value = [double]"3.893770667202106E-06"
$ws.Range("T16").Value = [double]"3.893770667202106E-06"
$ws.Range("E17").Value = [double]"1"
$ws.Range("F17").Value = [double]"0.3333333333333333"
$ws.Range("G17").Value = [double]"0.014556"
$ws.Range("H17").Value = [double]"0.043668"
$ws.Range("I17").Value = [double]"7.413387281768795E-05"
$ws.Range("J17").Value = [double]"7.413387281768795E-05"
$ws.Range("M17").Value = [double]"47.79503400000001"
$ws.Range("N17").Value = [double]"143.385102"
$ws.Range("O17").Value = [double]"0.4637252825399149"
$ws.Range("P17").Value = [double]"0.4637252825399149"
$ws.Range("Q17").Value = [double]"0.6957045149040001"
$ws.Range("R17").Value = [double]"6.261340634136"
$ws.Range("S17").Value = [double]"3.437775111816046E-05"
$ws.Range("T17").Value = [double]"0.001015599880090667"
